# Roboflow Annotation Report 7/18/2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Table1 currently spans D4:J60 (last data row = 60). Grow the table by one
# row via the ListObject so the table/autofilter ref (and sheet dimension)
# are updated to D4:J61, matching the diff.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.ListRows.Add() | Out-Null

$newRow = 61

# New weekly entry, mirrors the prior row's counts with the new date.
$ws.Range("D$newRow").Value = "18/7/2026"
$ws.Range("E$newRow").Value = 406
$ws.Range("F$newRow").Value = 924
$ws.Range("G$newRow").Value = 0
$ws.Range("H$newRow").Value = 0
$ws.Range("I$newRow").Value = 1012
$ws.Range("J$newRow").Value = "N/A"

# Copy the formatting (styles/number formats/borders) from the row above
# so the new row picks up the same look as the rest of the table.
$ws.Range("D60:J60").Copy()
$ws.Range("D$($newRow):J$newRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item($newRow).RowHeight = $ws.Rows.Item(60).RowHeight

# Update the view so the new row region is visible, matching the
# recorded sheetView state after the edit (scrolled down, selection F68).
$excel.ActiveWindow.ScrollRow = 47
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F68").Select()
